$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.300186991691589
$ws.Range("B1").Value = 2.247406721115112
$ws.Range("D1").Value = 1.897592663764954
$ws.Range("E1").Value = 0.8990330696105957
